$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 69.89967366666667
$ws.Range("H2").Value = 209.699021
$ws.Range("I2").Value = 0.6608367681537789
$ws.Range("J2").Value = 0.660836768153779
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 11750.82207075968
$ws.Range("R2").Value = 105757.3986368371
$ws.Range("S2").Value = 0.1972057358980056
$ws.Range("T2").Value = 0.1972057358980057

$ws.Range("G3").Value = 69.89967366666667
$ws.Range("H3").Value = 209.699021
$ws.Range("I3").Value = 0.6608367681537789
$ws.Range("J3").Value = 0.660836768153779
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 11394.08267873176
$ws.Range("R3").Value = 102546.7441085859
$ws.Range("S3").Value = 0.191218830989988
$ws.Range("T3").Value = 0.1912188309899881

$ws.Range("G4").Value = 69.89967366666667
$ws.Range("H4").Value = 209.699021
$ws.Range("I4").Value = 0.6608367681537789
$ws.Range("J4").Value = 0.660836768153779
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 11602.89357777804
$ws.Range("R4").Value = 104426.0422000024
$ws.Range("S4").Value = 0.1947231566245676
$ws.Range("T4").Value = 0.1947231566245676

$ws.Range("G5").Value = 69.89967366666667
$ws.Range("H5").Value = 209.699021
$ws.Range("I5").Value = 0.6608367681537789
$ws.Range("J5").Value = 0.660836768153779
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 4629.22711790903
$ws.Range("R5").Value = 41663.04406118127
$ws.Range("S5").Value = 0.07768904464121762
$ws.Range("T5").Value = 0.07768904464121763

$ws.Range("I6").Value = 0.1661491941864736
$ws.Range("J6").Value = 0.1661491941864736
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 2954.420383629476
$ws.Range("R6").Value = 26589.78345266528
$ws.Range("S6").Value = 0.04958194774776744
$ws.Range("T6").Value = 0.04958194774776745

$ws.Range("I7").Value = 0.1661491941864736
$ws.Range("J7").Value = 0.1661491941864736
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.04807670549419674
$ws.Range("T7").Value = 0.04807670549419676

$ws.Range("I8").Value = 0.1661491941864736
$ws.Range("J8").Value = 0.1661491941864736
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 2917.227840643733
$ws.Range("R8").Value = 26255.0505657936
$ws.Range("S8").Value = 0.04895777160372791
$ws.Range("T8").Value = 0.04895777160372791

$ws.Range("I9").Value = 0.1661491941864736
$ws.Range("J9").Value = 0.1661491941864736
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 1163.891587775236
$ws.Range("R9").Value = 10475.02428997712
$ws.Range("S9").Value = 0.01953276934078153
$ws.Range("T9").Value = 0.01953276934078153

$ws.Range("G10").Value = 4.152730666666667
$ws.Range("H10").Value = 12.458192
$ws.Range("I10").Value = 0.0392602278210887
$ws.Range("J10").Value = 0.03926022782108871
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 698.1148353351715
$ws.Range("R10").Value = 6283.033518016544
$ws.Range("S10").Value = 0.01171596753100076
$ws.Range("T10").Value = 0.01171596753100076

$ws.Range("G11").Value = 4.152730666666667
$ws.Range("H11").Value = 12.458192
$ws.Range("I11").Value = 0.0392602278210887
$ws.Range("J11").Value = 0.03926022782108871
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 676.9209937108604
$ws.Range("R11").Value = 6092.288943397744
$ws.Range("S11").Value = 0.01136028627662893
$ws.Range("T11").Value = 0.01136028627662893

$ws.Range("G12").Value = 4.152730666666667
$ws.Range("H12").Value = 12.458192
$ws.Range("I12").Value = 0.0392602278210887
$ws.Range("J12").Value = 0.03926022782108871
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 689.3264224992533
$ws.Range("R12").Value = 6203.93780249328
$ws.Range("S12").Value = 0.0115684778140902
$ws.Range("T12").Value = 0.0115684778140902

$ws.Range("G13").Value = 4.152730666666667
$ws.Range("H13").Value = 12.458192
$ws.Range("I13").Value = 0.0392602278210887
$ws.Range("J13").Value = 0.03926022782108871
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 275.0217906192196
$ws.Range("R13").Value = 2475.196115572976
$ws.Range("S13").Value = 0.00461549619936881
$ws.Range("T13").Value = 0.004615496199368811

$ws.Range("G14").Value = 14.14774133333333
$ws.Range("H14").Value = 42.443224
$ws.Range("I14").Value = 0.1337538098386587
$ws.Range("J14").Value = 0.1337538098386588
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 2378.37435270333
$ws.Range("R14").Value = 21405.36917432997
$ws.Range("S14").Value = 0.03991457462647808
$ws.Range("T14").Value = 0.0399145746264781

$ws.Range("G15").Value = 14.14774133333333
$ws.Range("H15").Value = 42.443224
$ws.Range("I15").Value = 0.1337538098386587
$ws.Range("J15").Value = 0.1337538098386588
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 2306.170057932374
$ws.Range("R15").Value = 20755.53052139137
$ws.Range("S15").Value = 0.03870282101472569
$ws.Range("T15").Value = 0.0387028210147257

$ws.Range("G16").Value = 14.14774133333333
$ws.Range("H16").Value = 42.443224
$ws.Range("I16").Value = 0.1337538098386587
$ws.Range("J16").Value = 0.1337538098386588
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 2348.433525446907
$ws.Range("R16").Value = 21135.90172902216
$ws.Range("S16").Value = 0.03941209889865726
$ws.Range("T16").Value = 0.03941209889865727

$ws.Range("G17").Value = 14.14774133333333
$ws.Range("H17").Value = 42.443224
$ws.Range("I17").Value = 0.1337538098386587
$ws.Range("J17").Value = 0.1337538098386588
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 936.958706699386
$ws.Range("R17").Value = 8432.628360294473
$ws.Range("S17").Value = 0.01572431529879769
$ws.Range("T17").Value = 0.0157243152987977
